$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 169 is a brand-new row (the table previously ended at row 168).
# Copy the A/E column formatting (bold+border+center style, and the date
# number format) from the last existing row so the new row matches the
# established per-column formatting instead of picking up the default style.
$ws.Range("A168").Copy()
$ws.Range("A169").PasteSpecial(-4122)
$ws.Range("E168").Copy()
$ws.Range("E169").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# A new match (Ferencvarosi TC vs Diosgyori VTK, played 2024-04-10) is
# inserted as row 163. It sorts chronologically before the existing row 163
# match (2024-04-13), so every match previously on rows 163-168 shifts down
# by one row, to rows 164-169. We rewrite every affected row's data in place
# (column A keeps its original per-row sequence number, since it is just the
# row's running index and does not travel with the shifted match data).

$ws.Range("A163").Value2 = 161
$ws.Range("B163").Value2 = 7939469
$ws.Range("C163").Value2 = "Hungary NB I"
$ws.Range("D163").Value2 = "Hungary NB I"
$ws.Range("E163").Value2 = 45392.54166666666
$ws.Range("F163").Value2 = "Ferencvarosi TC"
$ws.Range("G163").Value2 = "Diosgyori VTK"
$ws.Range("H163").Value2 = 2
$ws.Range("I163").Value2 = 1
$ws.Range("J163").Value2 = "H"
$ws.Range("K163").Value2 = 1.3
$ws.Range("L163").Value2 = 5.5
$ws.Range("M163").Value2 = 9
$ws.Range("N163").Value2 = 1.25
$ws.Range("O163").Value2 = 6.5
$ws.Range("P163").Value2 = 9.5
$ws.Range("Q163").Value2 = -1.75
$ws.Range("R163").Value2 = 1.875
$ws.Range("S163").Value2 = 1.975
$ws.Range("T163").Value2 = 3.25
$ws.Range("U163").Value2 = 1.975
$ws.Range("V163").Value2 = 1.875
$ws.Range("W163").Value2 = 0.25
$ws.Range("X163").Value2 = -1
$ws.Range("Y163").Value2 = -1
$ws.Range("Z163").Value2 = -1
$ws.Range("AA163").Value2 = 0.9750000000000001
$ws.Range("AB163").Value2 = -0.5
$ws.Range("AC163").Value2 = 0.4375
$ws.Range("A164").Value2 = 162
$ws.Range("B164").Value2 = 6818359
$ws.Range("C164").Value2 = "Hungary NB I"
$ws.Range("D164").Value2 = "Hungary NB I"
$ws.Range("E164").Value2 = 45395.39583333334
$ws.Range("F164").Value2 = "Kisvarda FC"
$ws.Range("G164").Value2 = "Ujpest"
$ws.Range("K164").Value2 = 2.6
$ws.Range("L164").Value2 = 3.4
$ws.Range("M164").Value2 = 2.6
$ws.Range("N164").Value2 = 2.5
$ws.Range("O164").Value2 = 3.4
$ws.Range("P164").Value2 = 2.7
$ws.Range("Q164").Value2 = 0
$ws.Range("R164").Value2 = 1.8
$ws.Range("S164").Value2 = 2.05
$ws.Range("T164").Value2 = 2.5
$ws.Range("U164").Value2 = 1.975
$ws.Range("V164").Value2 = 1.875
$ws.Range("W164").Value2 = 0
$ws.Range("X164").Value2 = 0
$ws.Range("Y164").Value2 = 0
$ws.Range("Z164").Value2 = 0
$ws.Range("AA164").Value2 = 0
$ws.Range("A165").Value2 = 163
$ws.Range("B165").Value2 = 6818358
$ws.Range("C165").Value2 = "Hungary NB I"
$ws.Range("D165").Value2 = "Hungary NB I"
$ws.Range("E165").Value2 = 45395.5
$ws.Range("F165").Value2 = "Debreceni VSC"
$ws.Range("G165").Value2 = "Kecskemeti TE"
$ws.Range("K165").Value2 = 1.666
$ws.Range("L165").Value2 = 3.5
$ws.Range("M165").Value2 = 5.5
$ws.Range("N165").Value2 = 1.571
$ws.Range("O165").Value2 = 3.75
$ws.Range("P165").Value2 = 6.5
$ws.Range("Q165").Value2 = -1
$ws.Range("R165").Value2 = 2.025
$ws.Range("S165").Value2 = 1.825
$ws.Range("T165").Value2 = 2.5
$ws.Range("U165").Value2 = 1.925
$ws.Range("V165").Value2 = 1.925
$ws.Range("W165").Value2 = 0
$ws.Range("X165").Value2 = 0
$ws.Range("Y165").Value2 = 0
$ws.Range("Z165").Value2 = 0
$ws.Range("AA165").Value2 = 0
$ws.Range("A166").Value2 = 164
$ws.Range("B166").Value2 = 6818360
$ws.Range("C166").Value2 = "Hungary NB I"
$ws.Range("D166").Value2 = "Hungary NB I"
$ws.Range("E166").Value2 = 45395.60416666666
$ws.Range("F166").Value2 = "Paksi"
$ws.Range("G166").Value2 = "MOL Fehervar FC"
$ws.Range("K166").Value2 = 1.666
$ws.Range("L166").Value2 = 3.5
$ws.Range("M166").Value2 = 5.5
$ws.Range("N166").Value2 = 1.7
$ws.Range("O166").Value2 = 3.6
$ws.Range("P166").Value2 = 5.25
$ws.Range("Q166").Value2 = -0.75
$ws.Range("R166").Value2 = 1.925
$ws.Range("S166").Value2 = 1.925
$ws.Range("T166").Value2 = 2.75
$ws.Range("U166").Value2 = 1.9
$ws.Range("V166").Value2 = 1.95
$ws.Range("W166").Value2 = 0
$ws.Range("X166").Value2 = 0
$ws.Range("Y166").Value2 = 0
$ws.Range("Z166").Value2 = 0
$ws.Range("AA166").Value2 = 0
$ws.Range("A167").Value2 = 165
$ws.Range("B167").Value2 = 6820520
$ws.Range("C167").Value2 = "Hungary NB I"
$ws.Range("D167").Value2 = "Hungary NB I"
$ws.Range("E167").Value2 = 45396.375
$ws.Range("F167").Value2 = "MTK Budapest"
$ws.Range("G167").Value2 = "Diosgyori VTK"
$ws.Range("K167").Value2 = 1.727
$ws.Range("L167").Value2 = 3.4
$ws.Range("M167").Value2 = 4.5
$ws.Range("N167").Value2 = 2.25
$ws.Range("O167").Value2 = 3.2
$ws.Range("P167").Value2 = 2.8
$ws.Range("Q167").Value2 = -0.25
$ws.Range("R167").Value2 = 2.05
$ws.Range("S167").Value2 = 1.8
$ws.Range("T167").Value2 = 2.75
$ws.Range("U167").Value2 = 1.85
$ws.Range("V167").Value2 = 2
$ws.Range("W167").Value2 = 0
$ws.Range("X167").Value2 = 0
$ws.Range("Y167").Value2 = 0
$ws.Range("Z167").Value2 = 0
$ws.Range("AA167").Value2 = 0
$ws.Range("A168").Value2 = 166
$ws.Range("B168").Value2 = 6860904
$ws.Range("C168").Value2 = "Hungary NB I"
$ws.Range("D168").Value2 = "Hungary NB I"
$ws.Range("E168").Value2 = 45396.47916666666
$ws.Range("F168").Value2 = "Zalaegerszegi TE"
$ws.Range("G168").Value2 = "Ferencvarosi TC"
$ws.Range("K168").Value2 = 5.5
$ws.Range("L168").Value2 = 4
$ws.Range("M168").Value2 = 1.5
$ws.Range("N168").Value2 = 5.75
$ws.Range("O168").Value2 = 4
$ws.Range("P168").Value2 = 1.5
$ws.Range("Q168").Value2 = 1.25
$ws.Range("R168").Value2 = 1.85
$ws.Range("S168").Value2 = 2
$ws.Range("T168").Value2 = 3.25
$ws.Range("U168").Value2 = 2
$ws.Range("V168").Value2 = 1.85
$ws.Range("W168").Value2 = 0
$ws.Range("X168").Value2 = 0
$ws.Range("Y168").Value2 = 0
$ws.Range("Z168").Value2 = 0
$ws.Range("AA168").Value2 = 0
$ws.Range("A169").Value2 = 167
$ws.Range("B169").Value2 = 6818361
$ws.Range("C169").Value2 = "Hungary NB I"
$ws.Range("D169").Value2 = "Hungary NB I"
$ws.Range("E169").Value2 = 45396.57291666666
$ws.Range("F169").Value2 = "Mezokovesd Zsory"
$ws.Range("G169").Value2 = "Puskas Academy"
$ws.Range("K169").Value2 = 5.5
$ws.Range("L169").Value2 = 3.5
$ws.Range("M169").Value2 = 1.666
$ws.Range("N169").Value2 = 5.5
$ws.Range("O169").Value2 = 3.6
$ws.Range("P169").Value2 = 1.65
$ws.Range("Q169").Value2 = 0.75
$ws.Range("R169").Value2 = 2
$ws.Range("S169").Value2 = 1.85
$ws.Range("T169").Value2 = 2.5
$ws.Range("U169").Value2 = 1.9
$ws.Range("V169").Value2 = 1.95
$ws.Range("W169").Value2 = 0
$ws.Range("X169").Value2 = 0
$ws.Range("Y169").Value2 = 0
$ws.Range("Z169").Value2 = 0
$ws.Range("AA169").Value2 = 0
